$wb = $excel.ActiveWorkbook

$hoUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d47472660c14f89d3bc9c8c0b32d44401898fae/e2e/ac048824-46ca-4c8f-a70c-1b15639240f9.md"
$hbUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8d47472660c14f89d3bc9c8c0b32d44401898fae/e2e/f2bfd839-ba30-4254-a8da-68a9fd3cc98c.md"

$linkColor = 15570276  # OLE (BGR) form of RGB(100,149,237) i.e. #FF6495ED "Cornflower Blue"

# -------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3)
# -------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

# -------------------------------------------------------------------
# 2. Overview sheet: widen the zh-cn / de-de columns (E, F)
# -------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.16666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.16666666666667

# -------------------------------------------------------------------
# 3. zh-cn sheet: fill in handback info for row 2 and row 3,
#    widen columns C, I, J
# -------------------------------------------------------------------
$wsZh.Columns.Item(3).ColumnWidth = 29.16666666666667
$wsZh.Columns.Item(9).ColumnWidth = 39.16666666666667
$wsZh.Columns.Item(10).ColumnWidth = 39.16666666666667

# Row 2 (ac048824-...)
$wsZh.Range("I2").Value = "ac048824-46ca-4c8f-a70c-1b15639240f9.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $hoUrl, [Type]::Missing, [Type]::Missing, "ac048824-46ca-4c8f-a70c-1b15639240f9.md")
$wsZh.Range("I2").Font.Color = $linkColor
$wsZh.Range("I2").Font.Underline = $true
$wsZh.Range("J2").Value = "ac048824-46ca-4c8f-a70c-1b15639240f9.f38a17eeac4cd81adcec9854be841540e8662918.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-19 12:46:35"

# Row 3 (f2bfd839-...)
$wsZh.Range("I3").Value = "f2bfd839-ba30-4254-a8da-68a9fd3cc98c.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $hbUrl, [Type]::Missing, [Type]::Missing, "f2bfd839-ba30-4254-a8da-68a9fd3cc98c.md")
$wsZh.Range("I3").Font.Color = $linkColor
$wsZh.Range("I3").Font.Underline = $true
$wsZh.Range("J3").Value = "f2bfd839-ba30-4254-a8da-68a9fd3cc98c.bf4abb0f080743712a7579d0beef3be1ec0ee763.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-19 12:46:35"

# -------------------------------------------------------------------
# 4. de-de sheet: fill in handback info for row 2 and row 3,
#    widen columns C, I, J
# -------------------------------------------------------------------
$wsDe.Columns.Item(3).ColumnWidth = 29.16666666666667
$wsDe.Columns.Item(9).ColumnWidth = 39.16666666666667
$wsDe.Columns.Item(10).ColumnWidth = 39.16666666666667

# Row 2 (ac048824-...)
$wsDe.Range("I2").Value = "ac048824-46ca-4c8f-a70c-1b15639240f9.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $hoUrl, [Type]::Missing, [Type]::Missing, "ac048824-46ca-4c8f-a70c-1b15639240f9.md")
$wsDe.Range("I2").Font.Color = $linkColor
$wsDe.Range("I2").Font.Underline = $true
$wsDe.Range("J2").Value = "ac048824-46ca-4c8f-a70c-1b15639240f9.f38a17eeac4cd81adcec9854be841540e8662918.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-19 12:46:42"

# Row 3 (f2bfd839-...)
$wsDe.Range("I3").Value = "f2bfd839-ba30-4254-a8da-68a9fd3cc98c.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $hbUrl, [Type]::Missing, [Type]::Missing, "f2bfd839-ba30-4254-a8da-68a9fd3cc98c.md")
$wsDe.Range("I3").Font.Color = $linkColor
$wsDe.Range("I3").Font.Underline = $true
$wsDe.Range("J3").Value = "f2bfd839-ba30-4254-a8da-68a9fd3cc98c.bf4abb0f080743712a7579d0beef3be1ec0ee763.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-19 12:46:42"
